# Generate Report for Handback
#
# The "Ready for handoff" status (shared by the Overview summary sheet and
# the per-language detail sheets) is replaced with "Handback transform
# failed", and the previously-empty "Error Detail" column (P) on the
# per-language sheets is populated with a diagnostic message for the
# 8d299c6b... file row. The Error Detail column is also widened so the
# message is readable.

$wb = $excel.ActiveWorkbook

# --- Update status text everywhere "Ready for handoff" was shown ---------
$newStatus = "Handback transform failed"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $newStatus

# --- Populate the Error Detail column (P) for the failed handback row ----
$wsZhCn.Range("P3").Value = "Handback file name: bzv3bjcf.3qz is different with handoff file name: 8d299c6b-10dc-45a9-8974-b95c2e657fb8.24cfe2742b00caf4265833273e01c7d10713b015.zh-cn."
$wsDeDe.Range("P3").Value = "Handback file name: bzv3bjcf.3qz is different with handoff file name: 8d299c6b-10dc-45a9-8974-b95c2e657fb8.24cfe2742b00caf4265833273e01c7d10713b015.de-de."

# --- Widen the Error Detail column so the new message is visible ---------
# ColumnWidth is in "characters"; Excel pads this by ~5/6 of a character
# when it round-trips through the stored OOXML <col width="..."/> value, so
# we back that padding out here to land on an effective width of 40.
$targetColWidth = 40 - (5 / 6)
$wsZhCn.Columns.Item(16).ColumnWidth = $targetColWidth
$wsDeDe.Columns.Item(16).ColumnWidth = $targetColWidth
